# updated legacy GSC export data
# Remove the 2025-10-17 row (incomplete/duplicate data point) from the
# "Chart" sheet. Deleting the entire row shifts every subsequent row up
# by one, which matches the upstream export after it was re-pulled.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
